# Edit script: add a "QuadracticCaliper(ms)" column (E) to the first table
# (A1:D8) on Sheet1, copying the header + values that already exist in the
# second small table at A30:B37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column E (same shared string as B30: "QuadracticCaliper(ms)")
$ws.Range("E1").Value = "QuadracticCaliper(ms)"

# Values for E2:E8 (copy of B31:B37)
$ws.Range("E2").Value = 0.91
$ws.Range("E3").Value = 0.83
$ws.Range("E4").Value = 4.88
$ws.Range("E5").Value = 20.19
$ws.Range("E6").Value = 146.44
$ws.Range("E7").Value = 616.52
$ws.Range("E8").Value = 3419.62

# Widen the new column to match the others roughly (target stored width ~27.33)
$ws.Range("E1").ColumnWidth = 26.5

# Update the selection to match the committed state
$ws.Activate()
$ws.Range("A21:B28").Select()
